$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing values in rows 2-5 ---
$ws.Range("B2").Value = 51.4
$ws.Range("C2").Value = 24.417

$ws.Range("B3").Value = 54.8
$ws.Range("C3").Value = 24.108

$ws.Range("B4").Value = 50.8
$ws.Range("C4").Value = 24.5

$ws.Range("B5").Value = 54.5
$ws.Range("C5").Value = 23.945

# --- Append new rows 6-9 ---
# Column A holds numeric-looking labels that must be stored as TEXT
# (matching the style/type of existing A2:A5 cells), so we force text
# via NumberFormat on a scratch cell, then restore the correct look
# (border/font/alignment) by copying the format from A5.
$newRows = @(
    @{ Row = 6; Label = "4"; Win = 55.7; Len = 24.47 },
    @{ Row = 7; Label = "5"; Win = 55.9; Len = 24.002 },
    @{ Row = 8; Label = "6"; Win = 52.3; Len = 24.369 },
    @{ Row = 9; Label = "7"; Win = 55;   Len = 24.221 }
)

foreach ($item in $newRows) {
    $r = $item.Row

    # Force the label into the scratch cell as text, then move the
    # (text-typed) value into the destination cell.
    $ws.Range("Z100").NumberFormat = "@"
    $ws.Range("Z100").Value = $item.Label
    $ws.Range("Z100").Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4163)  # xlPasteValues

    # Apply the same look (font/border/alignment) as the other label cells.
    $ws.Range("A5").Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)  # xlPasteFormats

    $ws.Cells.Item($r, 2).Value = $item.Win
    $ws.Cells.Item($r, 3).Value = $item.Len
}

# Remove the scratch cell entirely so it doesn't linger in the sheet.
$ws.Range("Z100").Delete(-4162)  # xlShiftUp
